$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The previous "latest" row (25) loses its date-only highlight format and
# reverts to the regular datetime number format used by all earlier rows.
$ws.Range("Y25").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Append the new row of bunker price data (row 26).
$rowValues = @{
    "A"  = 550
    "B"  = 462
    "C"  = 435
    "D"  = 535
    "E"  = 506
    "F"  = 514
    "G"  = 461
    "H"  = 555
    "I"  = 500
    "J"  = 435
    "K"  = 573
    "L"  = 465
    "M"  = 460
    "N"  = 491
    "O"  = 550
    "P"  = 465
    "Q"  = 620
    "R"  = 485
    "S"  = 461
    "T"  = 465
    "U"  = 621
    "V"  = 535
    "W"  = 590
    "X"  = 475
    "Y"  = 45756
    "Z"  = 848
    "AA" = 557
    "AB" = 535.5
    "AC" = 494
    "AD" = 543
    "AE" = 507
    "AF" = 509
    "AG" = 733
    "AH" = 455
    "AI" = 723
    "AJ" = 461
    "AK" = 472
    "AL" = 545
    "AM" = 536
    "AN" = 472
    "AO" = 519
    "AP" = 539
    "AQ" = 567
    "AR" = 550
    "AS" = 640
    "AT" = 639
    "AU" = 487
    "AV" = 462
}

foreach ($col in $rowValues.Keys) {
    $ws.Range($col + "26").Value = $rowValues[$col]
}

# New last row gets the date-only "latest entry" highlight format.
$ws.Range("Y26").NumberFormat = "YYYY-MM-DD"
